# Auditoría Configuración.xlsx - add new checklist items (project status per
# iteration, updated matrix/risk items, corrective-action template, etc.) to
# the "Documentos a revisar" list and retarget a couple of existing labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new / changed label text -------------------------------------------------
$estatus1 = "* Estatus del Proyecto en la 1era iteración"
$estatus2 = "* Estatus del Proyecto en la 2da iteración"
$matrizTrazabilidad = "* Matriz de trazabilidad"
$riesgos = "* Riesgos"
$componentes = "* Componentes del sistema"
$flujoVistas = "* Flujo de vistas"
$arquitecturaInicial = "* Arquitectura inicial"
$estimacionFP = "* Estimación de function points"
$plantillaAcciones = "* Plantilla de acciones correctivas"
$estatusGeneral = "* Estatus del Proyecto general"

# --- row 7 ---------------------------------------------------------------------
$ws.Range("H7").Value = "* Definición del proyecto"
$ws.Range("I7").Value = "Guía de estilos y estándares"
$ws.Range("J7").Value = $estatus1
$ws.Range("K7").Value = $estatus1

# --- row 8 ---------------------------------------------------------------------
$ws.Range("H8").Value = "*Propuesta de Proyecto"
$ws.Range("I8").Value = $estatus1
$ws.Range("K8").Value = $estatus2

# --- row 9 ---------------------------------------------------------------------
$ws.Range("H9").Value = $matrizTrazabilidad
$ws.Range("I9").Value = $estatus2

# --- row 10 ----------------------------------------------------------------------
$ws.Range("H10").Value = $riesgos

# --- row 11 ----------------------------------------------------------------------
$ws.Range("H11").Value = $componentes

# --- row 12 ----------------------------------------------------------------------
$ws.Range("H12").Value = $flujoVistas

# --- row 13 (new content, previously held the "* Arquitectura inicial" text) -----
$ws.Range("H13").Value = $arquitecturaInicial

# --- row 14 (Tarea/Encargado/Fecha límite header unaffected; H14 changes) --------
$ws.Range("H14").Value = $estimacionFP

# --- row 15 (brand new item) ------------------------------------------------------
$ws.Range("H15").Value = $plantillaAcciones

# --- rows 16-18 are brand new ----------------------------------------------------
$ws.Range("H16").Value = $estatus1
$ws.Range("H17").Value = $estatus2
$ws.Range("H18").Value = $estatusGeneral

# --- formatting: the whole "Documentos a revisar" column (H) and any populated
# I/J/K cells wrap their text, matching the rest of the checklist column -------
# (applied per contiguous block - the COM shim only honours WrapText on the
# first area of a multi-area union range)
$ws.Range("H7:H18").WrapText = $true
$ws.Range("I7:I9").WrapText = $true
$ws.Range("J7").WrapText = $true
$ws.Range("K7:K8").WrapText = $true

# Re-assert the centred/wrapped formatting on the "Actividad" column so it
# keeps round-tripping correctly alongside the rest of the sheet.
$ws.Range("A7:A11").WrapText = $true
$ws.Range("A7:A11").HorizontalAlignment = -4108

# --- row heights: rows that now hold wrapped 2-line captions ---------------------
$ws.Rows.Item(7).RowHeight = 46.25
$ws.Rows.Item(9).RowHeight = 46.25
$ws.Rows.Item(10).RowHeight = 35.05
$ws.Rows.Item(11).RowHeight = 46.25
$ws.Rows.Item(15).RowHeight = 23.85
$ws.Rows.Item(16).RowHeight = 23.85
$ws.Rows.Item(17).RowHeight = 23.85

# --- selection, matching the saved cursor position in the source workbook -------
$ws.Range("I11").Select() | Out-Null
